$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to stay text (some look like numbers/dates to Excel's
# auto-detection); write the value under a Text format, then restore the default
# "Normal" style so no extra formatting is left behind on the cell.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '72.238.25'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.633.82'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.10'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.91'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.44%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.632.59'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.114.25'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '72.162.87'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('E17').Value = '  -1.94%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.640.36'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.36'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.01%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.07'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '373.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('E22').Value = '  -2.00%  '
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '70.81'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.21'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.53'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.70%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.768.83'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0947'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.92'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.17%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '494.73'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.71%  '
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.06'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.114'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.23%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.84'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -6.77%  '
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.325'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '39.03'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '151.83'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('E50').Value = '  -2.56%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.600'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.24%  '
